# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) value or $null to leave unchanged, new Volume(1h) (E) value
$updates = @(
    @(2,  "71.753.37",  "  +4.01%  "),
    @(3,  "3.687.01",   "  +8.40%  "),
    @(4,  $null,        "  +0.00%  "),
    @(5,  "589.70",     "  +1.76%  "),
    @(6,  "181.17",     "  +2.26%  "),
    @(7,  "3.679.50",   "  +8.29%  "),
    @(8,  $null,        "  +4.71%  "),
    @(9,  $null,        "  +0.07%  "),
    @(10, "0.202",      "  +2.57%  "),
    @(11, $null,        "  +5.59%  "),
    @(12, "50.07",      "  +4.35%  "),
    @(13, $null,        "  +2.66%  "),
    @(14, "4.281.14",   "  +8.54%  "),
    @(15, "686.73",     "  +1.51%  "),
    @(16, "9.05",       "  +5.52%  "),
    @(17, "3.691.55",   "  +8.40%  "),
    @(18, "71.812.03",  "  +3.94%  "),
    @(19, $null,        "  +2.40%  "),
    @(20, "18.18",      "  +2.64%  "),
    @(21, "11.70",      "  +4.02%  "),
    @(22, $null,        "  +4.25%  "),
    @(23, "6.35",       "  +18.68%  "),
    @(24, "17.87",      "  +5.48%  "),
    @(25, "104.29",     "  +3.78%  "),
    @(26, "4.04",       "  +4.45%  "),
    @(27, "2.85",       "  +6.64%  "),
    @(28, "10.21",      "  +6.23%  "),
    @(29, "35.44",      "  +6.24%  "),
    @(30, "9.29",       "  +6.81%  "),
    @(31, "7.39",       "  +8.25%  "),
    @(32, "4.27",       "  +15.85%  "),
    @(33, $null,        "  +3.37%  "),
    @(34, "567.32",     "  +3.28%  "),
    @(35, $null,        "  +4.95%  "),
    @(36, "59.59",      "  +2.78%  "),
    @(37, "3.779.54",   "  +5.06%  "),
    @(38, $null,        "  -0.06%  "),
    @(39, "0.146",      "  +4.03%  "),
    @(40, $null,        "  +6.66%  "),
    @(41, "35.67",      "  +2.30%  "),
    @(42, $null,        "  +6.60%  "),
    @(43, $null,        "  +10.68%  "),
    @(44, "2.82",       "  +5.76%  "),
    @(45, "0.354",      "  +6.58%  "),
    @(46, "2.92",       "  +10.15%  "),
    @(47, "3.38",       "  +0.19%  "),
    @(48, $null,        "  +4.64%  "),
    @(49, $null,        "  +3.14%  "),
    @(50, "0.998",      "  -0.33%  "),
    @(51, "135.07",     "  +3.44%  ")
)

foreach ($entry in $updates) {
    $row = $entry[0]
    $price = $entry[1]
    $volume = $entry[2]

    if ($null -ne $price) {
        # The Price column holds plain text (e.g. "71.753.37" or "589.70").
        # A bare assignment lets Excel's COM layer auto-coerce single-decimal
        # looking text into a real number, which would change the cell's
        # stored type. Prefix with an apostrophe to force text entry, then
        # restore the default "Normal" style so the cell's formatting stays
        # identical to how it started (no visible/explicit style index).
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $price
        $cell.Style = "Normal"
    }

    $ws.Cells.Item($row, 5).Value = $volume
}
